$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alt3")

# Row 11 -> splits B3
$ws.Range("D11").Formula2 = '=TEXTSPLIT(B3,{"(",")","[","]","{","}"})'
$ws.Range("E11:F11").Style = "Normal"

# Row 12 -> splits B4
$ws.Range("D12").Formula2 = '=TEXTSPLIT(B4,{"(",")","[","]","{","}"})'
$ws.Range("E12:F12").Style = "Normal"

# Row 13 -> splits B5
$ws.Range("D13").Formula2 = '=TEXTSPLIT(B5,{"(",")","[","]","{","}"})'
$ws.Range("E13:F13").Style = "Normal"

# Row 14 -> splits B6
$ws.Range("D14").Formula2 = '=TEXTSPLIT(B6,{"(",")","[","]","{","}"})'
$ws.Range("E14:F14").Style = "Normal"

# Row 15 -> splits B7
$ws.Range("D15").Formula2 = '=TEXTSPLIT(B7,{"(",")","[","]","{","}"})'
$ws.Range("E15:F15").Style = "Normal"

$ws.Range("G8").Select() | Out-Null
